$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: now "Rewiring of light point..." item ---
$ws.Range("A8").Value = ""
$ws.Range("C8").Value = 37
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.0"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.00"

# --- Row 9: now "Rewiring of 3/5 pin..." item ---
$ws.Range("A9").Value = ""
$ws.Range("C9").Value = 62
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.0"
$ws.Range("E9").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.00"

# --- Row 10: now "P. point / On board" item (formerly row 8 content, qty/amount updated) ---
$ws.Range("A10").Value = "P. point"
$ws.Range("C10").Value = 91
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6"
$ws.Range("E10").Value = "On board"
$ws.Range("F10").Value = 136
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "12376.00"

# --- Row 11: now "Each / P & F ISI marked..." item ---
$ws.Range("A11").Value = "Each"
$ws.Range("C11").Value = 78
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.0"
$ws.Range("E11").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 50
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3900.00"

# --- Row 12: now "Each / 6 A to 32 A rating" item ---
$ws.Range("A12").Value = "Each"
$ws.Range("C12").Value = 3
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "30"
$ws.Range("E12").Value = ' 6 A to 32 A rating'
$ws.Range("F12").Value = 187
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "561.00"

# --- Totals rows updated to reflect new Grand Total (16837.00) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "16837.00"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "16837.00"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "16837.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "16837.00"

